$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp (08:50 -> 09:20)
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 09:20"

# Refresh country case-count data. Some countries overtook neighbours in the
# "Casos totales" ranking, so besides the numeric refresh the country name
# shown in a few rows also changes (the table stays sorted by column B desc).

$ws.Cells.Item(13, 1).Value = "Turquia"
$ws.Cells.Item(13, 2).Value = 13531
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(13, 4).Value = 243
$ws.Cells.Item(13, 5).Value = 13074
$ws.Cells.Item(13, 6).Value = 847
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 214

$ws.Cells.Item(16, 1).Value = "Austria"
$ws.Cells.Item(16, 2).Value = 10298
$ws.Cells.Item(16, 3).Value = 118
$ws.Cells.Item(16, 4).Value = 1095
$ws.Cells.Item(16, 5).Value = 9075
$ws.Cells.Item(16, 6).Value = 198
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = 128

$ws.Cells.Item(21, 1).Value = "Israel"
$ws.Cells.Item(21, 2).Value = 5591
$ws.Cells.Item(21, 3).Value = 233
$ws.Cells.Item(21, 4).Value = 226
$ws.Cells.Item(21, 5).Value = 5344
$ws.Cells.Item(21, 6).Value = 97
$ws.Cells.Item(21, 7).Value = 1
$ws.Cells.Item(21, 8).Value = 21

$ws.Cells.Item(25, 1).Value = "Chequia"
$ws.Cells.Item(25, 2).Value = 3330
$ws.Cells.Item(25, 3).Value = 22
$ws.Cells.Item(25, 4).Value = 45
$ws.Cells.Item(25, 5).Value = 3253
$ws.Cells.Item(25, 6).Value = 70
$ws.Cells.Item(25, 7).Value = 1
$ws.Cells.Item(25, 8).Value = 32

$ws.Cells.Item(38, 1).Value = "Tailandia"
$ws.Cells.Item(38, 2).Value = 1771
$ws.Cells.Item(38, 3).Value = 120
$ws.Cells.Item(38, 4).Value = 416
$ws.Cells.Item(38, 5).Value = 1343
$ws.Cells.Item(38, 6).Value = 23
$ws.Cells.Item(38, 7).Value = 2
$ws.Cells.Item(38, 8).Value = 12

$ws.Cells.Item(64, 1).Value = "Ucrania"
$ws.Cells.Item(64, 2).Value = 669
$ws.Cells.Item(64, 3).Value = 24
$ws.Cells.Item(64, 4).Value = 10
$ws.Cells.Item(64, 5).Value = 642
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 17

$ws.Cells.Item(65, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(65, 2).Value = 664
$ws.Cells.Item(65, 3).Value = 0
$ws.Cells.Item(65, 4).Value = 61
$ws.Cells.Item(65, 5).Value = 597
$ws.Cells.Item(65, 6).Value = 2
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 6

$ws.Cells.Item(66, 1).Value = "Marruecos"
$ws.Cells.Item(66, 2).Value = 638
$ws.Cells.Item(66, 3).Value = 21
$ws.Cells.Item(66, 4).Value = 24
$ws.Cells.Item(66, 5).Value = 578
$ws.Cells.Item(66, 6).Value = 1
$ws.Cells.Item(66, 7).Value = 0
$ws.Cells.Item(66, 8).Value = 36

$ws.Cells.Item(68, 1).Value = "Armenia"
$ws.Cells.Item(68, 2).Value = 571
$ws.Cells.Item(68, 3).Value = 39
$ws.Cells.Item(68, 4).Value = 31
$ws.Cells.Item(68, 5).Value = 537
$ws.Cells.Item(68, 6).Value = 30
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 3

$ws.Cells.Item(69, 1).Value = "Barein"
$ws.Cells.Item(69, 2).Value = 567
$ws.Cells.Item(69, 3).Value = 0
$ws.Cells.Item(69, 4).Value = 316
$ws.Cells.Item(69, 5).Value = 247
$ws.Cells.Item(69, 6).Value = 2
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 4

$ws.Cells.Item(73, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(73, 2).Value = 424
$ws.Cells.Item(73, 3).Value = 4
$ws.Cells.Item(73, 4).Value = 17
$ws.Cells.Item(73, 5).Value = 394
$ws.Cells.Item(73, 6).Value = 1
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 13

$ws.Cells.Item(83, 1).Value = "Taiwan"
$ws.Cells.Item(83, 2).Value = 329
$ws.Cells.Item(83, 3).Value = 7
$ws.Cells.Item(83, 4).Value = 45
$ws.Cells.Item(83, 5).Value = 279
$ws.Cells.Item(83, 6).Value = 0
$ws.Cells.Item(83, 7).Value = 0
$ws.Cells.Item(83, 8).Value = 5

$ws.Cells.Item(94, 1).Value = "Oman"
$ws.Cells.Item(94, 2).Value = 210
$ws.Cells.Item(94, 3).Value = 18
$ws.Cells.Item(94, 4).Value = 34
$ws.Cells.Item(94, 5).Value = 175
$ws.Cells.Item(94, 6).Value = 3
$ws.Cells.Item(94, 7).Value = 0
$ws.Cells.Item(94, 8).Value = 1

$ws.Cells.Item(95, 1).Value = "Afganistan"
$ws.Cells.Item(95, 2).Value = 196
$ws.Cells.Item(95, 3).Value = 22
$ws.Cells.Item(95, 4).Value = 5
$ws.Cells.Item(95, 5).Value = 187
$ws.Cells.Item(95, 6).Value = 0
$ws.Cells.Item(95, 7).Value = 0
$ws.Cells.Item(95, 8).Value = 4

$ws.Cells.Item(100, 1).Value = "Islas Feroe"
$ws.Cells.Item(100, 2).Value = 173
$ws.Cells.Item(100, 3).Value = 4
$ws.Cells.Item(100, 4).Value = 75
$ws.Cells.Item(100, 5).Value = 98
$ws.Cells.Item(100, 6).Value = 1
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = 0

$ws.Cells.Item(101, 1).Value = "Honduras"
$ws.Cells.Item(101, 2).Value = 172
$ws.Cells.Item(101, 3).Value = 31
$ws.Cells.Item(101, 4).Value = 3
$ws.Cells.Item(101, 5).Value = 159
$ws.Cells.Item(101, 6).Value = 4
$ws.Cells.Item(101, 7).Value = 3
$ws.Cells.Item(101, 8).Value = 10

$ws.Cells.Item(102, 1).Value = "Malta"
$ws.Cells.Item(102, 2).Value = 169
$ws.Cells.Item(102, 3).Value = 0
$ws.Cells.Item(102, 4).Value = 2
$ws.Cells.Item(102, 5).Value = 167
$ws.Cells.Item(102, 6).Value = 2
$ws.Cells.Item(102, 7).Value = 0
$ws.Cells.Item(102, 8).Value = 0

$ws.Cells.Item(111, 1).Value = "Montenegro"
$ws.Cells.Item(111, 2).Value = 120
$ws.Cells.Item(111, 3).Value = 11
$ws.Cells.Item(111, 4).Value = 0
$ws.Cells.Item(111, 5).Value = 118
$ws.Cells.Item(111, 6).Value = 4
$ws.Cells.Item(111, 7).Value = 0
$ws.Cells.Item(111, 8).Value = 2

$ws.Cells.Item(112, 1).Value = "Estado de Palestina"
$ws.Cells.Item(112, 2).Value = 119
$ws.Cells.Item(112, 3).Value = 0
$ws.Cells.Item(112, 4).Value = 18
$ws.Cells.Item(112, 5).Value = 100
$ws.Cells.Item(112, 6).Value = 0
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 1

$ws.Cells.Item(113, 1).Value = "Bolivia"
$ws.Cells.Item(113, 2).Value = 115
$ws.Cells.Item(113, 3).Value = 8
$ws.Cells.Item(113, 4).Value = 1
$ws.Cells.Item(113, 5).Value = 107
$ws.Cells.Item(113, 6).Value = 3
$ws.Cells.Item(113, 7).Value = 1
$ws.Cells.Item(113, 8).Value = 7

$ws.Cells.Item(114, 1).Value = "Georgia"
$ws.Cells.Item(114, 2).Value = 115
$ws.Cells.Item(114, 3).Value = 5
$ws.Cells.Item(114, 4).Value = 22
$ws.Cells.Item(114, 5).Value = 93
$ws.Cells.Item(114, 6).Value = 6
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 0

$ws.Cells.Item(115, 1).Value = "Guadalupe"
$ws.Cells.Item(115, 2).Value = 114
$ws.Cells.Item(115, 3).Value = 0
$ws.Cells.Item(115, 4).Value = 22
$ws.Cells.Item(115, 5).Value = 88
$ws.Cells.Item(115, 6).Value = 14
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 4

$ws.Cells.Item(116, 1).Value = "Kirguistan"
$ws.Cells.Item(116, 2).Value = 111
$ws.Cells.Item(116, 3).Value = 4
$ws.Cells.Item(116, 4).Value = 3
$ws.Cells.Item(116, 5).Value = 108
$ws.Cells.Item(116, 6).Value = 3
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 0
